$wb = $excel.ActiveWorkbook

# --- Sheet "Presentie": fill in the "Dag 2" (column D) attendance, mirroring
#     "Dag 1" (column C), except where column C says "ziek gemeld" (sick),
#     which becomes "-" for day 2. ---
$ws3 = $wb.Worksheets.Item("Presentie")

$dag2 = @{
    2  = "X"
    3  = "X"
    4  = "X"
    5  = "-"
    6  = "X"
    7  = "X"
    8  = "-"
    9  = "X"
    10 = "X"
    11 = "-"
    12 = "X"
    13 = "X"
}

foreach ($row in $dag2.Keys) {
    $cellC = $ws3.Cells.Item($row, 3)
    $cellD = $ws3.Cells.Item($row, 4)
    $cellD.Value = $dag2[$row]
    # Mirror Dag 1's cell formatting (incl. quote-prefixed "-" cells) onto Dag 2.
    $cellC.Copy()
    $cellD.PasteSpecial(-4122)
}

# Update the remembered selection on the Presentie sheet.
[void]$ws3.Range("D14").Select()

# --- Sheet "Inventarisatie": re-apply left alignment on B1:E1 (header row) ---
$ws1 = $wb.Worksheets.Item("Inventarisatie")
$ws1.Range("B1:E1").HorizontalAlignment = -4131
